# Sunil Narine .xlsx — add ownTeam/oppTeam columns, reorder rows by date,
# and append 3 newly-scraped matches (IPL / Kolkata Knight Riders data refresh).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet grows from 9 columns (A:I) / 7 rows to 11 columns (A:K) / 10 rows,
# and two new columns (ownTeam, oppTeam) are inserted before "batsman", so it
# is simplest to clear the old layout and rewrite the full table fresh.
$ws.Range("A1:K10").ClearContents()

# totalRuns/totalBalls/total4s/total6s/sr (columns G:K) are textual numbers
# (e.g. "100.00", "0.00") in the source data, so format them as Text first
# to stop Excel from normalising them into real numbers and dropping the
# trailing zeros / decimal formatting.
$ws.Range("G1:K10").NumberFormat = "@"

$table = @(
    @("venue", "date", "result", "ownTeam", "oppTeam", "batsman", "totalRuns", "totalBalls", "total4s", "total6s", "sr"),
    @(" Abu Dhabi", " October 07 2020", "KKR won by 10 runs", "Kolkata Knight Riders", "Chennai Super Kings", "Sunil Narine ", "17", "9", "1", "1", "188.88"),
    @(" Dubai (DSC)", " September 30 2020", "KKR won by 37 runs", "Kolkata Knight Riders", "Rajasthan Royals", "Sunil Narine ", "15", "14", "2", "1", "107.14"),
    @(" Dubai (DSC)", " November 01 2020", "KKR won by 60 runs", "Kolkata Knight Riders", "Rajasthan Royals", "Sunil Narine ", "0", "2", "0", "0", "0.00"),
    @(" Abu Dhabi", " September 26 2020", "KKR won by 7 wickets (with 12 balls remaining)", "Kolkata Knight Riders", "Sunrisers Hyderabad", "Sunil Narine ", "0", "2", "0", "0", "0.00"),
    @(" Sharjah", " October 03 2020", "Capitals won by 18 runs", "Kolkata Knight Riders", "Delhi Capitals", "Sunil Narine ", "3", "5", "0", "0", "60.00"),
    @(" Dubai (DSC)", " October 29 2020", "Super Kings won by 6 wickets", "Kolkata Knight Riders", "Chennai Super Kings", "Sunil Narine ", "7", "7", "0", "1", "100.00"),
    @(" Sharjah", " October 26 2020", "Kings XI won by 8 wickets (with 7 balls remaining)", "Kolkata Knight Riders", "Kings XI Punjab", "Sunil Narine ", "6", "4", "1", "0", "150.00"),
    @(" Abu Dhabi", " October 24 2020", "KKR won by 59 runs", "Kolkata Knight Riders", "Delhi Capitals", "Sunil Narine ", "64", "32", "6", "4", "200.00"),
    @(" Abu Dhabi", " September 23 2020", "Mumbai won by 49 runs", "Kolkata Knight Riders", "Mumbai Indians", "Sunil Narine ", "9", "10", "0", "1", "90.00")
)

for ($r = 0; $r -lt $table.Length; $r++) {
    $row = $table[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}
